$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cell B11 to hold the text value "1" (was "R40").
# A leading apostrophe forces Excel to store this as literal text instead
# of auto-converting the numeric-looking "1" into a Number.
$cell = $ws.Range("B11")
$cell.Value = "'1"
